# shield_driver_v0p3 BOM: swap the 2x2 male header for a 4-position
# receptacle (female) connector, per "male->female header in BOM".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Description (D2) and part number (E2) for the P0-P3 connector row.
$ws.Range("D2").Value = "4 Position Receptacle Connector 0.100"" (2.54mm) Through Hole Gold"
$ws.Range("E2").Value = "A26452-ND"

# New unit price for the replacement part; J2 (=H2*B2) and J15
# (=SUM(J2:J12)) recalculate automatically.
$ws.Range("H2").Value = 1.71

# Reflect the new active cell left behind by the edit.
$null = $ws.Range("A2").Select()
